$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New filename entries added to column A for rows 3-7 (shared strings 18-22)
$ws.Range("A3").Value = "SCRIPT/D73P11A/us0203.ssb"
$ws.Range("A4").Value = "SCRIPT/D73P11A/us0302.ssb"
$ws.Range("A5").Value = "SCRIPT/D73P11A/us0403.ssb"
$ws.Range("A6").Value = "SCRIPT/D73P11A/us2002.ssb"
$ws.Range("A7").Value = "SCRIPT/D73P11A/us2102.ssb"

# Rows 3-7 grow to accommodate the new wrapped text (43.2pt)
$ws.Range("A3:A7").RowHeight = 43.2

# Update the active selection to match the author's last-touched cell
[void]$ws.Range("E5").Select()

Write-Output "done"
